# Rtery failed testcases and IAnnotation Transformer to remove annotations from @Test
#
# - RUNMANAGER!C2 ("execute" for "loginlogoutTest"): yes -> no
# - DATA!B5 ("execute" for newtest/admin row): no -> yes
# - DATA sheet active selection moves to C5

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# Stop executing the loginlogoutTest row on RUNMANAGER.
$wsRunManager.Range("C2").Value = "no"

# Start executing the newtest/admin row on DATA.
$wsData.Range("B5").Value = "yes"

# Update the active selection on the DATA sheet.
$wsData.Range("C5").Select()
